$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Swap the match data (columns F:V) between rows 95 and 96 ---
# Columns A-E (index, pais, torneio, temporada, data_partida) stay the same.
$cols = 6..22   # F=6 ... V=22

$row95 = @{}
$row96 = @{}
foreach ($c in $cols) {
    $row95[$c] = $ws.Cells.Item(95, $c).Value()
    $row96[$c] = $ws.Cells.Item(96, $c).Value()
}
foreach ($c in $cols) {
    $ws.Cells.Item(95, $c).Value = $row96[$c]
    $ws.Cells.Item(96, $c).Value = $row95[$c]
}

# --- 2) Append new row 114 (Liberia vs Saprissa) ---
# Copy formatting (styles/number formats) from row 113 first, restricted to A:V
# so we don't touch/format the whole 16384-column row.
$ws.Range("A113:V113").Copy()
$ws.Range("A114:V114").PasteSpecial(-4122)  # xlPasteFormats

$ws.Cells.Item(114, 1).Value = 113
$ws.Cells.Item(114, 2).Value = "costa-rica"
$ws.Cells.Item(114, 3).Value = "primera-division"
$ws.Cells.Item(114, 4).Value = "2023-2024"
$ws.Cells.Item(114, 5).Value = 45240.125
$ws.Cells.Item(114, 6).Value = "Liberia"
$ws.Cells.Item(114, 7).Value = 1
$ws.Cells.Item(114, 8).Value = "Saprissa"
$ws.Cells.Item(114, 9).Value = 2
$ws.Cells.Item(114, 10).Value = 4.68
$ws.Cells.Item(114, 11).Value = "05/11/2023 18:13"
$ws.Cells.Item(114, 12).Value = 4.13
$ws.Cells.Item(114, 13).Value = "10/11/2023 02:58"
$ws.Cells.Item(114, 14).Value = 3.97
$ws.Cells.Item(114, 15).Value = "05/11/2023 18:13"
$ws.Cells.Item(114, 16).Value = 4.02
$ws.Cells.Item(114, 17).Value = "10/11/2023 02:58"
$ws.Cells.Item(114, 18).Value = 1.69
$ws.Cells.Item(114, 19).Value = "05/11/2023 18:13"
$ws.Cells.Item(114, 20).Value = 1.79
$ws.Cells.Item(114, 21).Value = "10/11/2023 02:58"
$ws.Cells.Item(114, 22).Value = "https://www.betexplorer.com/football/costa-rica/primera-division/liberia-saprissa/YVJEBkt5/"

Write-Host "Applied edits: swapped rows 95/96 and appended row 114."
